$wb = $excel.ActiveWorkbook

# --- NOESY sheet: clear the synced/imported peak-list values (A2:I13), ---
# --- keeping formatting, and update the selection to A2:K13 ---
$wsNOESY = $wb.Worksheets.Item("NOESY")
$wsNOESY.Activate()
$wsNOESY.Range("A2:I13").ClearContents()
$wsNOESY.Range("A2:K13").Select()

# --- Make "molecule" the active / selected tab (was "COSY") ---
$wsMolecule = $wb.Worksheets.Item("molecule")
$wsMolecule.Activate()
